# إضافة حدث جديد في Card10
# Fill in the previously-blank "nan" placeholder cells on row 22, then
# append a brand-new event row (23) with the same "nan" placeholder pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

# --- Row 22: the existing last row was missing its "nan" placeholders ---
# (columns B..K and M get the literal text "nan"; L, N, O already had values)
$row22Cols = @("B","C","D","E","F","G","H","I","J","K","M")
foreach ($col in $row22Cols) {
    $ws.Range("$col" + "22").Value = "nan"
}

# --- Row 23: new service event row ---
# leading apostrophe forces the numeric-looking "10" to be stored as text,
# matching every other "card" cell in column A on this sheet
$ws.Range("A23").Value = "'10"
$ws.Range("L23").Value = "14\4\2025"
$ws.Range("N23").Value = "تم تغييرالجرائد الاماميه (1_2_4_5_7_8) ومعايرها"
$ws.Range("O23").Value = "الخبير"

$wb.Save()
